$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (X1 crystal): part number + JLCPCB part change
$ws.Range("A4").Value = "CO32D6-125.000-33GDTSTL"
$ws.Range("D4").Value = "C5119019"

# Row 15 (was U2, now U5): part, designator, footprint, JLCPCB part change
$ws.Range("A15").Value = "MAX1853EXT+T"
$ws.Range("B15").Value = "U5"
$ws.Range("C15").Value = "SOT-363-6"
$ws.Range("D15").Value = "C143384"

# Row 16 (R11): fill in previously empty JLCPCB part#
$ws.Range("D16").Value = "C23204"

# Row 17 (R12): fill in previously empty JLCPCB part#
$ws.Range("D17").Value = "C22787"
